$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 1.02
$ws.Range("F3").Value = 3.3
$ws.Range("I3").Value = 2.44
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 3.8
$ws.Range("L3").Value = 1.35
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.4
$ws.Range("O3").Value = 1.28
$ws.Range("P3").Value = 1.82
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.31
$ws.Range("S3").Value = 3.55
$ws.Range("T3").Value = 1.76
$ws.Range("U3").Value = 2.04
$ws.Range("V3").Value = 1.69
$ws.Range("W3").Value = 1.34
$ws.Range("X3").Value = 16
$ws.Range("Y3").Value = 11.5
$ws.Range("Z3").Value = 17.5
$ws.Range("AA3").Value = 980
$ws.Range("AB3").Value = 16
$ws.Range("AC3").Value = 9.4
$ws.Range("AD3").Value = 13.5
$ws.Range("AE3").Value = 980
$ws.Range("AF3").Value = 980
$ws.Range("AG3").Value = 18.5
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 980
$ws.Range("AJ3").Value = 85
$ws.Range("AK3").Value = 55
$ws.Range("AL3").Value = 65
$ws.Range("AM3").Value = 130
$ws.Range("AN3").Value = 60
$ws.Range("AO3").Value = 980
$ws.Range("G5").Value = 2.82
$ws.Range("V5").Value = 1.48
$ws.Range("W5").Value = 1.55
$ws.Range("G6").Value = 1.31
$ws.Range("J6").Value = 5.9
$ws.Range("K6").Value = 8.199999999999999
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 1.52
$ws.Range("W6").Value = 4.2
$ws.Range("F7").Value = 9.800000000000001
$ws.Range("G7").Value = 16.5
$ws.Range("L7").Value = 1.32
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.28
$ws.Range("R7").Value = 1.35
$ws.Range("S7").Value = 3
$ws.Range("T7").Value = 2.28
$ws.Range("U7").Value = 1.6
$ws.Range("V7").Value = 3.2
$ws.Range("X7").Value = 980
$ws.Range("Y7").Value = 8.199999999999999
$ws.Range("Z7").Value = 8.6
$ws.Range("AA7").Value = 980
$ws.Range("AB7").Value = 980
$ws.Range("AC7").Value = 980
$ws.Range("AD7").Value = 980
$ws.Range("AE7").Value = 980
$ws.Range("AG7").Value = 60
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 65
$ws.Range("AO7").Value = 8.6
$ws.Range("F8").Value = 3.15
$ws.Range("G8").Value = 3.9
$ws.Range("I8").Value = 2.62
$ws.Range("N8").Value = 3.15
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.75
$ws.Range("Q8").Value = 2.08
$ws.Range("S8").Value = 3.75
$ws.Range("T8").Value = 1.81
$ws.Range("V8").Value = 1.62
$ws.Range("W8").Value = 1.34
$ws.Range("Z8").Value = 18.5
$ws.Range("AA8").Value = 42
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 9.199999999999999
$ws.Range("AD8").Value = 14
$ws.Range("AG8").Value = 17.5
$ws.Range("AJ8").Value = 75
$ws.Range("AL8").Value = 70
$ws.Range("AO8").Value = 30
$ws.Range("P9").Value = 2.16
$ws.Range("Q9").Value = 1.63
$ws.Range("AO9").Value = 27
$ws.Range("G10").Value = 1.92
$ws.Range("J10").Value = 3.5
$ws.Range("T10").Value = 1.82
$ws.Range("W10").Value = 2.08
$ws.Range("F11").Value = 1.76
$ws.Range("I11").Value = 6
$ws.Range("Q11").Value = 1.88
$ws.Range("S11").Value = 3.25
$ws.Range("F12").Value = 3.35
$ws.Range("G12").Value = 4.2
$ws.Range("H12").Value = 2.06
$ws.Range("I12").Value = 2.48
$ws.Range("K12").Value = 4.3
$ws.Range("Q12").Value = 1.58
$ws.Range("R12").Value = 1.4
$ws.Range("S12").Value = 2.4
$ws.Range("V12").Value = 1.67
$ws.Range("W12").Value = 1.31
$ws.Range("F13").Value = 1.3
$ws.Range("G13").Value = 1.34
$ws.Range("H13").Value = 9.6
$ws.Range("I13").Value = 11.5
$ws.Range("K13").Value = 7.2
$ws.Range("L13").Value = 1.19
$ws.Range("Q13").Value = 1.41
$ws.Range("R13").Value = 1.83
$ws.Range("S13").Value = 2
$ws.Range("U13").Value = 2.04
$ws.Range("V13").Value = 1.09
$ws.Range("W13").Value = 3.9
$ws.Range("Y13").Value = 48
$ws.Range("Z13").Value = 130
$ws.Range("AA13").Value = 380
$ws.Range("AC13").Value = 16.5
$ws.Range("AE13").Value = 160
$ws.Range("AH13").Value = 27
$ws.Range("AI13").Value = 120
$ws.Range("AM13").Value = 130
$ws.Range("AN13").Value = 3.95
$ws.Range("AO13").Value = 150
$ws.Range("G14").Value = 3.2
$ws.Range("N14").Value = 4.6
$ws.Range("R14").Value = 1.59
$ws.Range("S14").Value = 2.38
$ws.Range("W14").Value = 1.46
$ws.Range("F15").Value = 3.2
$ws.Range("G15").Value = 4.1
$ws.Range("I15").Value = 2.32
$ws.Range("J15").Value = 3.85
$ws.Range("K15").Value = 6
$ws.Range("L15").Value = 1.2
$ws.Range("N15").Value = 3.3
$ws.Range("O15").Value = 1.09
$ws.Range("P15").Value = 3.3
$ws.Range("Q15").Value = 1.29
$ws.Range("R15").Value = 1.84
$ws.Range("S15").Value = 1.79
$ws.Range("T15").Value = 1.34
$ws.Range("V15").Value = 1.75
$ws.Range("W15").Value = 1.32
$ws.Range("X15").Value = 65
$ws.Range("AF15").Value = 50
$ws.Range("AJ15").Value = 85
$ws.Range("AN15").Value = 21
$ws.Range("AO15").Value = 9.800000000000001
$ws.Range("H17").Value = 3.25
$ws.Range("F18").Value = 1.51
$ws.Range("G18").Value = 1.53
$ws.Range("H18").Value = 8.6
$ws.Range("I18").Value = 8.800000000000001
$ws.Range("J18").Value = 4.4
$ws.Range("K18").Value = 4.5
$ws.Range("L18").Value = 1.42
$ws.Range("N18").Value = 3.45
$ws.Range("P18").Value = 1.83
$ws.Range("T18").Value = 2.32
$ws.Range("U18").Value = 1.7
$ws.Range("W18").Value = 2.9
$ws.Range("AA18").Value = 390
$ws.Range("AF18").Value = 7.6
$ws.Range("AJ18").Value = 12.5
$ws.Range("AL18").Value = 50
$ws.Range("AN18").Value = 9.4
$ws.Range("AO18").Value = 290
$ws.Range("F19").Value = 2.4
$ws.Range("G19").Value = 2.72
$ws.Range("I19").Value = 4.4
$ws.Range("K19").Value = 3.25
$ws.Range("U19").Value = 1.69
$ws.Range("V19").Value = 1.32
$ws.Range("W19").Value = 1.58
$ws.Range("F20").Value = 2.08
$ws.Range("G20").Value = 2.1
$ws.Range("H20").Value = 4.3
$ws.Range("I20").Value = 4.4
$ws.Range("R20").Value = 1.3
$ws.Range("V20").Value = 1.29
$ws.Range("W20").Value = 1.91
$ws.Range("Y20").Value = 13.5
$ws.Range("AJ20").Value = 24
$ws.Range("F21").Value = 7
$ws.Range("I21").Value = 1.6
$ws.Range("J21").Value = 4.1
$ws.Range("U21").Value = 1.91
$ws.Range("V21").Value = 2.66
$ws.Range("J22").Value = 3.65
$ws.Range("L22").Value = 1.33
$ws.Range("V22").Value = 1.3
$ws.Range("W22").Value = 1.87

Write-Output "Applied 199 cell updates"
